$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New workout rows to append (ExerciseId, DateId, Exercise Date, Exercise Week,
# Exercise Month, Exercise Year, Exercise Day, Exercise Name, Weight, Sets, Reps)
$rows = @(
    @(289, 35, 43096, 52, "December", 2017, "Wednesday", "Bench Press", 100, 5, 5),
    @(290, 35, 43096, 52, "December", 2017, "Wednesday", "Overhead Press", 55, 5, 5),
    @(291, 35, 43096, 52, "December", 2017, "Wednesday", "Barbell Row", 90, 5, 5),
    @(292, 35, 43096, 52, "December", 2017, "Wednesday", "Pull-Ups", 105, 5, 5),
    @(293, 35, 43096, 52, "December", 2017, "Wednesday", "Bicycles", 0, 3, 10),
    @(294, 35, 43096, 52, "December", 2017, "Wednesday", "Leg Raises", 0, 3, 10),
    @(295, 35, 43096, 52, "December", 2017, "Wednesday", "Russian Twists", 10, 3, 10),
    @(296, 35, 43096, 52, "December", 2017, "Wednesday", "V-up crunches with medicine ball", 8, 3, 10),
    @(297, 35, 43097, 52, "December", 2017, "Thursday", "Shoulder Press", 26, 4, 8),
    @(298, 35, 43097, 52, "December", 2017, "Thursday", "Shoulder Shrug", 26, 4, 8),
    @(299, 35, 43097, 52, "December", 2017, "Thursday", "Dumbell Chest Press", 26, 4, 8),
    @(300, 35, 43097, 52, "December", 2017, "Thursday", "One arm row (left)", 32, 4, 8),
    @(301, 35, 43097, 52, "December", 2017, "Thursday", "One arm row (right)", 32, 4, 8),
    @(302, 35, 43097, 52, "December", 2017, "Thursday", "Seated Row", 59, 4, 8),
    @(303, 35, 43097, 52, "December", 2017, "Thursday", "Lat Pull Down", 66, 4, 8),
    @(304, 35, 43097, 52, "December", 2017, "Thursday", "Left Crunch", 0, 4, 10),
    @(305, 35, 43097, 52, "December", 2017, "Thursday", "Right Crunch", 0, 4, 10),
    @(306, 35, 43097, 52, "December", 2017, "Thursday", "Sled Pushes", 30, 5, 4)
)

$startRow = 290
$endRow = $startRow + $rows.Count - 1

# Copy cell formatting (number formats / styles) from the last existing data
# row down onto the new rows before filling in values.
$ws.Range("A289:K289").Copy()
$ws.Range("A290:K$endRow").PasteSpecial(-4122)  # xlPasteFormats

$r = $startRow
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $r = $r + 1
}

$win = $excel.ActiveWindow
$win.ScrollRow = 277
[void]$ws.Range("C311").Select()
